$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1196
$ws.Range("I28").Value = 1234.5333
$ws.Range("J28").Value = 1113.4286
$ws.Range("K28").Value = 1234.5333
$ws.Range("L28").Value = 1113.4286
$ws.Range("M28").Value = -749.5333000000001
$ws.Range("N28").Value = -2083.4286
$ws.Range("H87").Value = 84333.336
$ws.Range("J87").Value = 84333.336
$ws.Range("L87").Value = 84333.336
$ws.Range("N87").Value = -86829.336
$ws.Range("H88").Value = 1397.9445
$ws.Range("I88").Value = 817
$ws.Range("J88").Value = 1862.7
$ws.Range("K88").Value = 817
$ws.Range("L88").Value = 1862.7
$ws.Range("M88").Value = -411
$ws.Range("N88").Value = -2674.7
$ws.Range("H90").Value = 84333.336
$ws.Range("J90").Value = 84333.336
$ws.Range("L90").Value = 253000.008
$ws.Range("N90").Value = -265480.008
$ws.Range("H91").Value = 1397.9445
$ws.Range("I91").Value = 817
$ws.Range("J91").Value = 1862.7
$ws.Range("K91").Value = 817
$ws.Range("L91").Value = 1862.7
$ws.Range("M91").Value = 587
$ws.Range("N91").Value = -4670.7
$ws.Range("H111").Value = 15691.111
$ws.Range("I111").Value = 16402.5
$ws.Range("K111").Value = 49207.5
$ws.Range("M111").Value = -46140.5
$ws.Range("H132").Value = 19728084
$ws.Range("I132").Value = 20898588
$ws.Range("K132").Value = 62695764
$ws.Range("M132").Value = -62693234
$ws.Range("H137").Value = 19752.848
$ws.Range("I137").Value = 15337.467
$ws.Range("J137").Value = 29214.38
$ws.Range("K137").Value = 46012.401
$ws.Range("L137").Value = 87643.14
$ws.Range("M137").Value = -43462.401
$ws.Range("N137").Value = -92743.14
$ws.Range("H138").Value = 4652.11
$ws.Range("I138").Value = 2629.3
$ws.Range("J138").Value = 4876.8667
$ws.Range("K138").Value = 7887.900000000001
$ws.Range("L138").Value = 14630.6001
$ws.Range("M138").Value = -2747.900000000001
$ws.Range("N138").Value = -24910.6001

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3064.5186
$ws.Range("I61").Value = 2817.8696
$ws.Range("J61").Value = 4482.75
$ws.Range("K61").Value = 2817.8696
$ws.Range("L61").Value = 4482.75
$ws.Range("M61").Value = -2605.8696
$ws.Range("N61").Value = -4906.75
$ws.Range("H74").Value = 32406.086
$ws.Range("I74").Value = 32034.781
$ws.Range("K74").Value = 32034.781
$ws.Range("M74").Value = -31160.781
$ws.Range("H77").Value = 32406.086
$ws.Range("I77").Value = 32034.781
$ws.Range("K77").Value = 160173.905
$ws.Range("M77").Value = -155805.905
$ws.Range("H132").Value = 4104.4644
$ws.Range("I132").Value = 3587.5908
$ws.Range("K132").Value = 10762.7724
$ws.Range("M132").Value = -8232.7724
$ws.Range("H136").Value = 3064.5186
$ws.Range("I136").Value = 2817.8696
$ws.Range("J136").Value = 4482.75
$ws.Range("K136").Value = 8453.6088
$ws.Range("L136").Value = 13448.25
$ws.Range("M136").Value = -5903.6088
$ws.Range("N136").Value = -18548.25

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 157
$ws.Range("I11").Value = 110.5
$ws.Range("J11").Value = 250
$ws.Range("K11").Value = 110.5
$ws.Range("L11").Value = 250
$ws.Range("M11").Value = 29.5
$ws.Range("N11").Value = -530
$ws.Range("H81").Value = 58364.43
$ws.Range("J81").Value = 58364.43
$ws.Range("L81").Value = 58364.43
$ws.Range("N81").Value = -60486.43
$ws.Range("H84").Value = 58364.43
$ws.Range("J84").Value = 58364.43
$ws.Range("L84").Value = 175093.29
$ws.Range("N84").Value = -185701.29
$ws.Range("H107").Value = 2926.0952
$ws.Range("I107").Value = 2449.9092
$ws.Range("J107").Value = 3449.9
$ws.Range("K107").Value = 2449.9092
$ws.Range("L107").Value = 3449.9
$ws.Range("M107").Value = -529.9092000000001
$ws.Range("N107").Value = -7289.9
$ws.Range("H134").Value = 2651.926
$ws.Range("I134").Value = 2209.5
$ws.Range("K134").Value = 6628.5
$ws.Range("M134").Value = -4093.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 382665
$ws.Range("J9").Value = 382665
$ws.Range("L9").Value = 382665
$ws.Range("N9").Value = -383001
$ws.Range("H31").Value = 35950.39
$ws.Range("I31").Value = 42439.05
$ws.Range("J31").Value = 12355.272
$ws.Range("K31").Value = 42439.05
$ws.Range("L31").Value = 12355.272
$ws.Range("M31").Value = -42144.05
$ws.Range("N31").Value = -12945.272
$ws.Range("H34").Value = 35950.39
$ws.Range("I34").Value = 42439.05
$ws.Range("J34").Value = 12355.272
$ws.Range("K34").Value = 42439.05
$ws.Range("L34").Value = 12355.272
$ws.Range("M34").Value = -42237.05
$ws.Range("N34").Value = -12759.272
$ws.Range("H105").Value = 3275
$ws.Range("I105").Value = 3275
$ws.Range("K105").Value = 3275
$ws.Range("M105").Value = -1528
$ws.Range("H107").Value = 1442
$ws.Range("I107").Value = 1457.7778
$ws.Range("K107").Value = 1457.7778
$ws.Range("M107").Value = 462.2221999999999
$ws.Range("H132").Value = 3402.0356
$ws.Range("I132").Value = 3268.7778
$ws.Range("K132").Value = 9806.3334
$ws.Range("M132").Value = -7276.3334

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 26991350
$ws.Range("I4").Value = 32091858
$ws.Range("K4").Value = 96275574
$ws.Range("M4").Value = -96275462
$ws.Range("H63").Value = 9000
$ws.Range("J63").Value = 9000
$ws.Range("L63").Value = 27000
$ws.Range("N63").Value = -28498
$ws.Range("H66").Value = 9000
$ws.Range("J66").Value = 9000
$ws.Range("L66").Value = 81000
$ws.Range("N66").Value = -88488
$ws.Range("H68").Value = 1517638.5
$ws.Range("J68").Value = 2176648.8
$ws.Range("L68").Value = 6529946.399999999
$ws.Range("N68").Value = -6531568.399999999
$ws.Range("H70").Value = 5341.3335
$ws.Range("I70").Value = 562
$ws.Range("J70").Value = 14900
$ws.Range("K70").Value = 1686
$ws.Range("L70").Value = 44700
$ws.Range("M70").Value = -1371
$ws.Range("N70").Value = -45330
$ws.Range("H71").Value = 1517638.5
$ws.Range("J71").Value = 2176648.8
$ws.Range("L71").Value = 19589839.2
$ws.Range("N71").Value = -19597951.2
$ws.Range("H73").Value = 5341.3335
$ws.Range("I73").Value = 562
$ws.Range("J73").Value = 14900
$ws.Range("K73").Value = 1686
$ws.Range("L73").Value = 44700
$ws.Range("M73").Value = -594
$ws.Range("N73").Value = -46884
$ws.Range("H107").Value = 848
$ws.Range("I107").Value = 481.25
$ws.Range("J107").Value = 1581.5
$ws.Range("K107").Value = 1443.75
$ws.Range("L107").Value = 4744.5
$ws.Range("M107").Value = 476.25
$ws.Range("N107").Value = -8584.5
$ws.Range("H112").Value = 11332.667
$ws.Range("J112").Value = 12999
$ws.Range("L112").Value = 38997
$ws.Range("N112").Value = -41213
$ws.Range("H113").Value = 303.75
$ws.Range("I113").Value = 313.05264
$ws.Range("K113").Value = 939.15792
$ws.Range("M113").Value = 1230.84208
$ws.Range("H122").Value = 1394.7142
$ws.Range("J122").Value = 2353.6
$ws.Range("L122").Value = 21182.4
$ws.Range("N122").Value = -26082.4
$ws.Range("H129").Value = 5502094
$ws.Range("I129").Value = 14143632
$ws.Range("J129").Value = 2934
$ws.Range("K129").Value = 42430896
$ws.Range("L129").Value = 8802
$ws.Range("M129").Value = -42425896
$ws.Range("N129").Value = -18802
$ws.Range("H131").Value = 18992.018
$ws.Range("J131").Value = 2579.1226
$ws.Range("L131").Value = 7737.3678
$ws.Range("N131").Value = -17817.3678

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1570
$ws.Range("I22").Value = 1570
$ws.Range("K22").Value = 1570
$ws.Range("M22").Value = -1275
$ws.Range("H27").Value = 1570
$ws.Range("I27").Value = 1570
$ws.Range("K27").Value = 1570
$ws.Range("M27").Value = -1463
$ws.Range("H136").Value = 5856.857
$ws.Range("I136").Value = 5856.857
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 17570.571
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -15020.571
$ws.Range("N136").ClearContents()
